$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.757.42'
$ws.Range("E2").Value = '  +2.66%  '
$ws.Range("D3").Value = '1.875.33'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '324.97'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '0.4594'
$ws.Range("E7").Value = '  -0.97%  '
$ws.Range("D8").Value = '0.3868'
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.07858'
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("D10").Value = '0.9886'
$ws.Range("E10").Value = '  +3.08%  '
$ws.Range("D11").Value = '21.80'
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '1.871.10'
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("D13").Value = '6.998'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").Value = '5.708'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '0.06942'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").Value = '88.45'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '0.00001003'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '16.79'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '28.761.18'
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").Value = '5.283'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '11.04'
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '2.099'
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("D25").Value = '2.078.91'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").Value = '152.89'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").Value = '19.24'
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").Value = '5.847'
$ws.Range("E28").Value = '  +3.48%  '
$ws.Range("D29").Value = '1.972'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").Value = '119.01'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").Value = '0.09320'
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").Value = '0.9193'
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").Value = '5.306'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("D34").Value = '1.338'
$ws.Range("E34").Value = '  +1.62%  '
$ws.Range("D35").Value = '3.323'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("D36").Value = '0.05774'
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").Value = '1.154'
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("D38").Value = '0.02077'
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("D39").Value = '7.686'
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = '0.5643'
$ws.Range("E40").Value = '  +1.05%  '
$ws.Range("D41").Value = '0.1787'
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("D42").Value = '9.859'
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '0.07217'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '11.75'
$ws.Range("E44").Value = '  +1.59%  '
$ws.Range("D45").Value = '0.5294'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").Value = '2.140'
$ws.Range("E46").Value = '  +2.41%  '
$ws.Range("D47").Value = '1.126'
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").Value = '1.830'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").Value = '113.49'
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("D50").Value = '2.411'
$ws.Range("E50").Value = '  +3.83%  '
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  +0.37%  '
